# Add team record (Wins/Losses/Ties) columns to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns AD, AE, AF ---
# Copy formatting (font/border/alignment) from an existing header cell (A1)
# so the new header cells share the same style as the rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-47): team record is constant for every player row ---
$wins = 107
$losses = 55
$ties = 0

for ($r = 2; $r -le 47; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins    # column AD = 30
    $ws.Cells.Item($r, 31).Value = $losses  # column AE = 31
    $ws.Cells.Item($r, 32).Value = $ties    # column AF = 32
}

Write-Output "Team record columns added"
